# Weekly update: insert 2 new observation rows for "Pepino ensalada" at
# Macroferia Regional de Talca, pushing the existing rows 533..599 down to
# 535..601 and growing the used range from A1:R599 to A1:R601.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 533 (shifts 533:599 -> 535:601).
$ws.Rows("533:534").Insert()

# --- New row 533: Región de Arica y Parinacota, $/caja 60 unidades ---
$ws.Cells.Item(533, 1).Value = 5
$ws.Cells.Item(533, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(533, 3).Value = "Maule"
$ws.Cells.Item(533, 4).Value = 44984
$ws.Cells.Item(533, 5).Value = 7
$ws.Cells.Item(533, 6).Value = 100112043
$ws.Cells.Item(533, 7).Value = "Pepino ensalada"
$ws.Cells.Item(533, 8).Value = "Sin especificar"
$ws.Cells.Item(533, 9).Value = "Primera"
$ws.Cells.Item(533, 10).Value = 400
$ws.Cells.Item(533, 11).Value = 5500
$ws.Cells.Item(533, 12).Value = 5500
$ws.Cells.Item(533, 13).Value = 5500
$ws.Cells.Item(533, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(533, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(533, 16).Value = 92
$ws.Cells.Item(533, 17).Value = 60
$ws.Cells.Item(533, 18).Value = "Hortaliza"

# --- New row 534: Región del Maule, $/caja 80 unidades ---
$ws.Cells.Item(534, 1).Value = 5
$ws.Cells.Item(534, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(534, 3).Value = "Maule"
$ws.Cells.Item(534, 4).Value = 44984
$ws.Cells.Item(534, 5).Value = 7
$ws.Cells.Item(534, 6).Value = 100112043
$ws.Cells.Item(534, 7).Value = "Pepino ensalada"
$ws.Cells.Item(534, 8).Value = "Sin especificar"
$ws.Cells.Item(534, 9).Value = "Primera"
$ws.Cells.Item(534, 10).Value = 400
$ws.Cells.Item(534, 11).Value = 8000
$ws.Cells.Item(534, 12).Value = 8000
$ws.Cells.Item(534, 13).Value = 8000
$ws.Cells.Item(534, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(534, 15).Value = "Región del Maule"
$ws.Cells.Item(534, 16).Value = 100
$ws.Cells.Item(534, 17).Value = 80
$ws.Cells.Item(534, 18).Value = "Hortaliza"
